$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Engineering" (A5) to "Process" and move it after "Logistics" (A4),
# shifting Production/Management/Store down and adding a thin border around
# the new "Process" cell.
$ws.Range("A5").Value = "Process"
$ws.Range("A6").Value = "Production"
$ws.Range("A7").Value = "Management"
$ws.Range("A8").Value = "Store"

# Apply a thin box border around A5 (the relocated "Process" cell).
$ws.Range("A5").Borders.LineStyle = 1   # xlContinuous
$ws.Range("A5").Borders.Weight = 2      # xlThin

# Update the selection to match the final state of the workbook.
$ws.Range("B11").Select()
